$d = $word.ActiveDocument

# List of (old, new) text replacements taken from the diff.
$replacements = @(
    @("2025-09-04 Thursday", "2025-09-05 Friday"),
    @("518÷5=103, 3", "610÷6=101, 4"),
    @("871÷3=290, 1", "232÷8=29, 0"),
    @("816÷4=204, 0", "743÷6=123, 5"),
    @("750÷7=107, 1", "989÷7=141, 2"),
    @("648÷4=162, 0", "895÷7=127, 6"),
    @("356÷8=44, 4", "625÷5=125, 0"),
    @("459÷6=76, 3", "663÷2=331, 1"),
    @("704÷7=100, 4", "640÷4=160, 0"),
    @("154÷5=30, 4", "135÷8=16, 7"),
    @("755÷8=94, 3", "910÷6=151, 4"),
    @("887÷5=177, 2", "140÷9=15, 5"),
    @("506÷4=126, 2", "242÷6=40, 2"),
    @("491÷9=54, 5", "653÷6=108, 5"),
    @("721÷8=90, 1", "951÷8=118, 7"),
    @("638÷5=127, 3", "826÷4=206, 2"),
    @("668÷4=167, 0", "318÷6=53, 0"),
    @("703÷3=234, 1", "230÷3=76, 2"),
    @("505÷7=72, 1", "102÷9=11, 3"),
    @("466÷5=93, 1", "695÷9=77, 2"),
    @("262÷3=87, 1", "793÷8=99, 1"),
    @("637÷6=106, 1", "863÷4=215, 3"),
    @("562÷2=281, 0", "732÷4=183, 0"),
    @("482÷3=160, 2", "505÷8=63, 1"),
    @("815÷3=271, 2", "567÷5=113, 2"),
    @("164÷9=18, 2", "495÷7=70, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
